$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 147; this shifts existing rows 147-155 down to 148-156
$ws.Rows.Item(147).Insert()

# Fill in the new row 147 with the boilerplate values copied from the (now shifted) row 148,
# plus the specific changed values from the diff.
$ws.Range("A147").Value = 10
$ws.Range("B147").Value = "Vega Modelo de Temuco"
$ws.Range("C147").Value = "La Araucanía"
$ws.Range("D147").Value = 45267
$ws.Range("E147").Value = 9
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100101
$ws.Range("H147").Value = "Berries"
$ws.Range("I147").Value = 100101001
$ws.Range("J147").Value = "Arándano (blue)"
$ws.Range("K147").Value = "Sin especificar"
$ws.Range("L147").Value = "Primera"
$ws.Range("M147").Value = 250
$ws.Range("N147").Value = 3300
$ws.Range("O147").Value = 3300
$ws.Range("P147").Value = 3300
$ws.Range("Q147").Value = "$/kilo"
$ws.Range("R147").Value = "Región del Maule"
$ws.Range("S147").Value = 3300
$ws.Range("T147").Value = 1
